$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999999999993094
$ws.Range("E2").Value = 0.9999999999993094

# Row 3
$ws.Range("D3").Value = 0.9911635519119191
$ws.Range("E3").Value = 0.9911635519119191

# Row 4
$ws.Range("D4").Value = 0.9999999846458427
$ws.Range("E4").Value = 0.9999999846458427

# Row 5
$ws.Range("D5").Value = 0.01633123592702319
$ws.Range("E5").Value = 0.01633123592702319

# Row 6
$ws.Range("D6").Value = [double]"6.511060816153663E-14"
$ws.Range("E6").Value = [double]"6.511060816153663E-14"

# Row 7
$ws.Range("D7").Value = 0.0002322839017102757
$ws.Range("E7").Value = 0.9997677160982897

# Row 8
$ws.Range("D8").Value = [double]"1.730902965311407E-05"
$ws.Range("E8").Value = 0.9999826909703469

# Row 9
$ws.Range("D9").Value = 0.9999999999654348
$ws.Range("E9").Value = [double]"3.456523955946977E-11"

# Row 10
$ws.Range("D10").Value = [double]"5.027881224770894E-09"
$ws.Range("E10").Value = 0.9999999949721188

# Row 11
$ws.Range("D11").Value = 0.9987628137152652
$ws.Range("E11").Value = 0.001237186284734793
$ws.Range("F11").Value = 8.917984008789062
